# edit.ps1
# Applies updated crypto price/volume values to sheet1 (rows 2-51, columns D and E)
# as described by the commit diff. All target cells hold plain text (inline strings)
# in the original workbook, so we force a Text number format before writing the
# value to prevent Excel from auto-converting numeric-looking strings (e.g. "512.18")
# into floating point numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ R = 2; C = 'D'; V = '57.752.16' },
    @{ R = 2; C = 'E'; V = '  -0.44%  ' },
    @{ R = 3; C = 'D'; V = '2.463.12' },
    @{ R = 3; C = 'E'; V = '  +0.42%  ' },
    @{ R = 4; C = 'E'; V = '  -0.02%  ' },
    @{ R = 5; C = 'D'; V = '512.18' },
    @{ R = 5; C = 'E'; V = '  -2.35%  ' },
    @{ R = 6; C = 'D'; V = '134.09' },
    @{ R = 6; C = 'E'; V = '  +3.20%  ' },
    @{ R = 7; C = 'E'; V = '  -0.22%  ' },
    @{ R = 8; C = 'D'; V = '0.559' },
    @{ R = 8; C = 'E'; V = '  -1.25%  ' },
    @{ R = 9; C = 'D'; V = '2.463.82' },
    @{ R = 9; C = 'E'; V = '  +0.20%  ' },
    @{ R = 10; C = 'D'; V = '0.0985' },
    @{ R = 10; C = 'E'; V = '  +0.85%  ' },
    @{ R = 11; C = 'E'; V = '  -0.44%  ' },
    @{ R = 12; C = 'E'; V = '  +0.88%  ' },
    @{ R = 13; C = 'D'; V = '4.65' },
    @{ R = 13; C = 'E'; V = '  -6.22%  ' },
    @{ R = 14; C = 'D'; V = '2.896.94' },
    @{ R = 14; C = 'E'; V = '  +0.22%  ' },
    @{ R = 15; C = 'D'; V = '57.772.97' },
    @{ R = 15; C = 'E'; V = '  -0.35%  ' },
    @{ R = 16; C = 'D'; V = '22.06' },
    @{ R = 16; C = 'E'; V = '  +2.22%  ' },
    @{ R = 17; C = 'E'; V = '  +1.65%  ' },
    @{ R = 18; C = 'D'; V = '2.465.61' },
    @{ R = 18; C = 'E'; V = '  +0.35%  ' },
    @{ R = 19; C = 'E'; V = '  +0.12%  ' },
    @{ R = 20; C = 'E'; V = '  +0.99%  ' },
    @{ R = 21; C = 'D'; V = '315.75' },
    @{ R = 21; C = 'E'; V = '  +1.32%  ' },
    @{ R = 22; C = 'D'; V = '6.51' },
    @{ R = 22; C = 'E'; V = '  +6.35%  ' },
    @{ R = 23; C = 'D'; V = '0.999' },
    @{ R = 23; C = 'E'; V = '  -0.05%  ' },
    @{ R = 24; C = 'E'; V = '  -1.76%  ' },
    @{ R = 25; C = 'D'; V = '65.35' },
    @{ R = 25; C = 'E'; V = '  +0.66%  ' },
    @{ R = 26; C = 'D'; V = '0.998' },
    @{ R = 26; C = 'E'; V = '  -0.23%  ' },
    @{ R = 27; C = 'E'; V = '  -0.04%  ' },
    @{ R = 28; C = 'D'; V = '0.384' },
    @{ R = 28; C = 'E'; V = '  -4.31%  ' },
    @{ R = 29; C = 'D'; V = '7.66' },
    @{ R = 29; C = 'E'; V = '  +5.63%  ' },
    @{ R = 30; C = 'D'; V = '172.40' },
    @{ R = 30; C = 'E'; V = '  -1.37%  ' },
    @{ R = 31; C = 'D'; V = '0.0₃0740' },
    @{ R = 31; C = 'E'; V = '  +0.68%  ' },
    @{ R = 32; C = 'E'; V = '  +0.58%  ' },
    @{ R = 33; C = 'D'; V = '6.17' },
    @{ R = 33; C = 'E'; V = '  +0.06%  ' },
    @{ R = 34; C = 'D'; V = '1.15' },
    @{ R = 34; C = 'E'; V = '  +1.10%  ' },
    @{ R = 35; C = 'E'; V = '  +0.04%  ' },
    @{ R = 36; C = 'D'; V = '0.995' },
    @{ R = 36; C = 'E'; V = '  -0.24%  ' },
    @{ R = 37; C = 'E'; V = '  +1.51%  ' },
    @{ R = 38; C = 'E'; V = '  +5.47%  ' },
    @{ R = 39; C = 'E'; V = '  +2.90%  ' },
    @{ R = 40; C = 'D'; V = '36.83' },
    @{ R = 40; C = 'E'; V = '  +1.34%  ' },
    @{ R = 41; C = 'E'; V = '  +1.94%  ' },
    @{ R = 42; C = 'D'; V = '0.809' },
    @{ R = 42; C = 'E'; V = '  +0.40%  ' },
    @{ R = 43; C = 'D'; V = '136.60' },
    @{ R = 43; C = 'E'; V = '  +9.66%  ' },
    @{ R = 44; C = 'E'; V = '  +1.16%  ' },
    @{ R = 45; C = 'D'; V = '4.99' },
    @{ R = 45; C = 'E'; V = '  +4.07%  ' },
    @{ R = 46; C = 'D'; V = '257.96' },
    @{ R = 46; C = 'E'; V = '  -0.07%  ' },
    @{ R = 47; C = 'D'; V = '0.579' },
    @{ R = 47; C = 'E'; V = '  -0.96%  ' },
    @{ R = 48; C = 'D'; V = '0.0923' },
    @{ R = 48; C = 'E'; V = '  +0.03%  ' },
    @{ R = 49; C = 'D'; V = '0.0496' },
    @{ R = 49; C = 'E'; V = '  +1.00%  ' },
    @{ R = 50; C = 'D'; V = '0.0216' },
    @{ R = 50; C = 'E'; V = '  +2.27%  ' },
    @{ R = 51; C = 'D'; V = '17.30' },
    @{ R = 51; C = 'E'; V = '  +1.45%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Cells.Item($u.R, $u.C)
    $cell.NumberFormat = "@"
    $cell.Value = $u.V
}
